$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the "2022-Q2" sheet (to inherit identical sheetPr/styles/column
#        formatting), placing the copy before it, then rename + overwrite its data
#        with the new 2022-Q3 figures. ---
$src = $wb.Worksheets.Item("2022-Q2")
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$src.Copy($beforeSheet)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q3"

# --- 2. Force text columns (B, D, E, F, G) to avoid numeric auto-conversion ---
$ws.Range("B2:B13").NumberFormat = "@"
$ws.Range("D2:G13").NumberFormat = "@"

# --- 3. Overwrite header row (values only; formatting already inherited) ---
$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# --- 4. Overwrite data rows (row 2..13) with the 2022-Q3 figures ---
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "003567"
$ws.Cells.Item(2,3).Value = "华夏行业景气混合"
$ws.Cells.Item(2,4).Value = "115.66"
$ws.Cells.Item(2,5).Value = "88.33"
$ws.Cells.Item(2,6).Value = "2.47"
$ws.Cells.Item(2,7).Value = "2.8568"
$ws.Cells.Item(2,8).Value = 5
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "006348"
$ws.Cells.Item(3,3).Value = "银华盛利混合A"
$ws.Cells.Item(3,4).Value = "12.76"
$ws.Cells.Item(3,5).Value = "86.66"
$ws.Cells.Item(3,6).Value = "3.18"
$ws.Cells.Item(3,7).Value = "0.4058"
$ws.Cells.Item(3,8).Value = 6
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "202019"
$ws.Cells.Item(4,3).Value = "南方策略优化混合"
$ws.Cells.Item(4,4).Value = "2.88"
$ws.Cells.Item(4,5).Value = "93.97"
$ws.Cells.Item(4,6).Value = "2.34"
$ws.Cells.Item(4,7).Value = "0.0674"
$ws.Cells.Item(4,8).Value = 6
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "001728"
$ws.Cells.Item(5,3).Value = "银华战略新兴灵活配置定期开放混合"
$ws.Cells.Item(5,4).Value = "1.39"
$ws.Cells.Item(5,5).Value = "97.07"
$ws.Cells.Item(5,6).Value = "3.32"
$ws.Cells.Item(5,7).Value = "0.0461"
$ws.Cells.Item(5,8).Value = 7
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "015684"
$ws.Cells.Item(6,3).Value = "银华盛利混合C"
$ws.Cells.Item(6,4).Value = "1.05"
$ws.Cells.Item(6,5).Value = "86.66"
$ws.Cells.Item(6,6).Value = "3.18"
$ws.Cells.Item(6,7).Value = "0.0334"
$ws.Cells.Item(6,8).Value = 6
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "002145"
$ws.Cells.Item(7,3).Value = "诺安景鑫灵活配置混合"
$ws.Cells.Item(7,4).Value = "0.50"
$ws.Cells.Item(7,5).Value = "77.19"
$ws.Cells.Item(7,6).Value = "3.87"
$ws.Cells.Item(7,7).Value = "0.0194"
$ws.Cells.Item(7,8).Value = 10
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "011231"
$ws.Cells.Item(8,3).Value = "光大保德信锦弘混合A"
$ws.Cells.Item(8,4).Value = "1.95"
$ws.Cells.Item(8,5).Value = "26.05"
$ws.Cells.Item(8,6).Value = "0.83"
$ws.Cells.Item(8,7).Value = "0.0162"
$ws.Cells.Item(8,8).Value = 3
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "620004"
$ws.Cells.Item(9,3).Value = "金元顺安价值增长混合"
$ws.Cells.Item(9,4).Value = "0.34"
$ws.Cells.Item(9,5).Value = "74.99"
$ws.Cells.Item(9,6).Value = "1.99"
$ws.Cells.Item(9,7).Value = "0.0068"
$ws.Cells.Item(9,8).Value = 6
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "011232"
$ws.Cells.Item(10,3).Value = "光大保德信锦弘混合C"
$ws.Cells.Item(10,4).Value = "0.82"
$ws.Cells.Item(10,5).Value = "26.05"
$ws.Cells.Item(10,6).Value = "0.83"
$ws.Cells.Item(10,7).Value = "0.0068"
$ws.Cells.Item(10,8).Value = 3
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "006157"
$ws.Cells.Item(11,3).Value = "财通量化核心优选混合"
$ws.Cells.Item(11,4).Value = "0.11"
$ws.Cells.Item(11,5).Value = "83.87"
$ws.Cells.Item(11,6).Value = "1.89"
$ws.Cells.Item(11,7).Value = "0.0021"
$ws.Cells.Item(11,8).Value = 2
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "011987"
$ws.Cells.Item(12,3).Value = "财通资管智选核心回报6个月持有期混合A"
$ws.Cells.Item(12,4).Value = "0.12"
$ws.Cells.Item(12,5).Value = "39.46"
$ws.Cells.Item(12,6).Value = "1.16"
$ws.Cells.Item(12,7).Value = "0.0014"
$ws.Cells.Item(12,8).Value = 9
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "011988"
$ws.Cells.Item(13,3).Value = "财通资管智选核心回报6个月持有期混合C"
$ws.Cells.Item(13,4).Value = "0.01"
$ws.Cells.Item(13,5).Value = "39.46"
$ws.Cells.Item(13,6).Value = "1.16"
$ws.Cells.Item(13,7).Value = "0.0001"
$ws.Cells.Item(13,8).Value = 9

# --- 5. Drop the temporary text-number-format residue (keeps values as text) ---
$ws.Range("B2:B13").ClearFormats()
$ws.Range("D2:G13").ClearFormats()

# Re-apply style "2" (bold/centered/bordered) to header + index column, since the
# ClearFormats() above only touched B/D:G columns, so A-column/header style (set by
# the sheet duplication) is untouched, but re-assert it defensively via copy/paste of
# formats from a sheet that is guaranteed to keep that exact style.
$styleSrc = $wb.Worksheets.Item("2022-Q1")
$styleSrc.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$ws.Range("A2:A13").PasteSpecial(-4122)

# --- 6. Update "总计" summary sheet: insert new row 2 for 2022-Q3 ---
$total = $wb.Worksheets.Item("总计")
$total.Cells.Item(2,1).EntireRow.Insert()
$total.Range("A2:D2").ClearFormats()
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 12
$total.Cells.Item(2,4).Value = 3.46

$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)

# Re-number the shifted-down rows' index column (A) so it stays 0,1,2,3
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3

# --- 7. Restore original active sheet/tab selection ---
$total.Activate()
